$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

$ws.Range("D2").Value = "58.125.33"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "2.361.78"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  -0.11%  "
$scratch.Formula = '="544.23"'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +2.67%  "
$scratch.Formula = '="136.17"'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +5.40%  "
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("E10").Value = "  +3.70%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("E12").Value = "  +1.18%  "
$scratch.Formula = '="24.03"'
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("D14").Value = "2.783.00"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "58.102.32"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").Value = "2.358.62"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("E18").Value = "  +3.75%  "
$scratch.Formula = '="333.08"'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("E20").Value = "  +2.47%  "
$scratch.Formula = '="6.79"'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("E28").Value = "  +2.46%  "
$scratch.Formula = '="172.33"'
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  +11.67%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("E34").Value = "  +0.05%  "
$scratch.Formula = '="4.25"'
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +6.98%  "
$ws.Range("E36").Value = "  +0.72%  "
$scratch.Formula = '="1.25"'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +3.80%  "
$ws.Range("E39").Value = "  +0.58%  "
$scratch.Formula = '="145.21"'
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -2.75%  "
$scratch.Formula = '="293.43"'
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +3.18%  "
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("E43").Value = "  +1.58%  "
$scratch.Formula = '="0.0949"'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +2.43%  "
$scratch.Formula = '="19.34"'
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("E48").Value = "  +3.07%  "
$scratch.Formula = '="17.56"'
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("E50").Value = "  +0.14%  "
$scratch.Formula = '="11.07"'
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.47%  "

$scratch.ClearContents()
$excel.CutCopyMode = $false
